# foaie_parcurs_B-151-VGT_iunie_2022_Alex_Bora.xlsx
# Update the "Km initiali" starting odometer reading, the daily trip log
# (km travelled / destination / reason) for the days that moved around,
# the monthly totals that derive from them, and the hand-off date in the
# signature line at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Starting odometer reading (Km initiali:)
$ws.Range("B12").Value = 43773

# Daily log rows (A=Ziua, B=Km_parcursi, C=Locul deplasarii, D=Observatii utilizator)
# Row 15 -> day 2
$ws.Range("B15").Value = 30
$ws.Range("C15").Value = "Acasa-Birou"
$ws.Range("D15").Value = " "

# Row 16 -> day 3
$ws.Range("B16").Value = 30
$ws.Range("C16").Value = "Acasa-Birou"
$ws.Range("D16").Value = " "

# Row 19 -> day 6
$ws.Range("B19").Value = 92
$ws.Range("C19").Value = "Cluj-Bontida"
$ws.Range("D19").Value = "Interes Serviciu"

# Row 21 -> day 8
$ws.Range("B21").Value = 421
$ws.Range("C21").Value = "Cluj-Satu-Mare"
$ws.Range("D21").Value = "Interes Serviciu"

# Row 22 -> day 9
$ws.Range("B22").Value = 257
$ws.Range("C22").Value = "Cluj-Bistrita"
$ws.Range("D22").Value = "Interes Serviciu"

# Row 23 -> day 10
$ws.Range("B23").Value = 30
$ws.Range("C23").Value = "Acasa-Birou"
$ws.Range("D23").Value = " "

# Row 27 -> day 14
$ws.Range("B27").Value = 101
$ws.Range("C27").Value = "Cluj-Dej"
$ws.Range("D27").Value = "Interes Serviciu"

# Row 33 -> day 20
$ws.Range("B33").Value = 356
$ws.Range("C33").Value = "Cluj-Baia-Mare"
$ws.Range("D33").Value = "Interes Serviciu"

# Row 34 -> day 21
$ws.Range("B34").Value = 85
$ws.Range("C34").Value = "Cluj-Apahida"
$ws.Range("D34").Value = "Interes Serviciu"

# Row 35 -> day 22
$ws.Range("B35").Value = 30
$ws.Range("C35").Value = "Acasa-Birou"
$ws.Range("D35").Value = " "

# Row 36 -> day 23
$ws.Range("B36").Value = 156
$ws.Range("C36").Value = "Cluj-Zalau"
$ws.Range("D36").Value = "Interes Serviciu"

# Row 37 -> day 24
$ws.Range("B37").Value = 30
$ws.Range("C37").Value = "Acasa-Birou"
$ws.Range("D37").Value = " "

# Row 40 -> day 27
$ws.Range("B40").Value = 257
$ws.Range("C40").Value = "Cluj-Bistrita"
$ws.Range("D40").Value = "Interes Serviciu"

# Row 42 -> day 29 (new destination string)
$ws.Range("B42").Value = 152
$ws.Range("C42").Value = "Cluj-Cmp. Turzii"
$ws.Range("D42").Value = "Interes Serviciu"

# Monthly totals (Km parcursi: / Total)
$ws.Range("B44").Value = 2712
$ws.Range("B45").Value = 46485

# Hand-off date in the signature line
$ws.Range("A55").Value = "Semnătură utilizator:`t`t`t  Data predarii: 01.07.2022"
